$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the cryptos list.
# NumberFormat is forced to text ("@") before assignment so that values
# which look numeric (e.g. "1.002", "1.000") are kept as literal text,
# matching the original inline-string cell content. The style is reset
# back to "Normal" afterwards so no stray cell formatting is introduced.
$cells = @(
    @{ Ref = "D2"; Value = "30.002.92" }
    @{ Ref = "E2"; Value = "  -0.46%  " }
    @{ Ref = "D3"; Value = "1.871.32" }
    @{ Ref = "D4"; Value = "1.002" }
    @{ Ref = "E4"; Value = "  +0.13%  " }
    @{ Ref = "D5"; Value = "319.50" }
    @{ Ref = "E5"; Value = "  -2.36%  " }
    @{ Ref = "D6"; Value = "1.001" }
    @{ Ref = "E6"; Value = "  +0.07%  " }
    @{ Ref = "D7"; Value = "0.5080" }
    @{ Ref = "E7"; Value = "  -1.85%  " }
    @{ Ref = "E8"; Value = "  -2.94%  " }
    @{ Ref = "D9"; Value = "0.08195" }
    @{ Ref = "E9"; Value = "  -3.36%  " }
    @{ Ref = "D10"; Value = "42.20" }
    @{ Ref = "E10"; Value = "  -1.56%  " }
    @{ Ref = "D11"; Value = "1.093" }
    @{ Ref = "E11"; Value = "  -3.12%  " }
    @{ Ref = "D12"; Value = "22.82" }
    @{ Ref = "E12"; Value = "  +2.83%  " }
    @{ Ref = "D13"; Value = "1.867.28" }
    @{ Ref = "E13"; Value = "  -3.05%  " }
    @{ Ref = "D14"; Value = "6.269" }
    @{ Ref = "E14"; Value = "  -1.52%  " }
    @{ Ref = "D15"; Value = "7.166" }
    @{ Ref = "E15"; Value = "  -3.03%  " }
    @{ Ref = "E16"; Value = "  +0.10%  " }
    @{ Ref = "D17"; Value = "91.91" }
    @{ Ref = "E17"; Value = "  -4.30%  " }
    @{ Ref = "D18"; Value = "0.00001083" }
    @{ Ref = "E18"; Value = "  -2.99%  " }
    @{ Ref = "D19"; Value = "0.06430" }
    @{ Ref = "E19"; Value = "  -4.65%  " }
    @{ Ref = "D20"; Value = "17.91" }
    @{ Ref = "E20"; Value = "  -2.03%  " }
    @{ Ref = "E21"; Value = "  +0.07%  " }
    @{ Ref = "D22"; Value = "29.990.16" }
    @{ Ref = "E22"; Value = "  -0.46%  " }
    @{ Ref = "D23"; Value = "5.809" }
    @{ Ref = "E23"; Value = "  -4.24%  " }
    @{ Ref = "E24"; Value = "  -1.91%  " }
    @{ Ref = "D25"; Value = "2.142" }
    @{ Ref = "E25"; Value = "  -2.61%  " }
    @{ Ref = "D26"; Value = "2.077.65" }
    @{ Ref = "E26"; Value = "  -3.24%  " }
    @{ Ref = "D27"; Value = "161.22" }
    @{ Ref = "E27"; Value = "  +0.27%  " }
    @{ Ref = "D28"; Value = "20.94" }
    @{ Ref = "E28"; Value = "  -1.69%  " }
    @{ Ref = "D29"; Value = "2.239" }
    @{ Ref = "E29"; Value = "  -9.13%  " }
    @{ Ref = "D30"; Value = "127.08" }
    @{ Ref = "E30"; Value = "  -1.76%  " }
    @{ Ref = "D31"; Value = "1.055" }
    @{ Ref = "D32"; Value = "0.1034" }
    @{ Ref = "E32"; Value = "  -2.47%  " }
    @{ Ref = "D33"; Value = "5.905" }
    @{ Ref = "E33"; Value = "  -3.17%  " }
    @{ Ref = "D34"; Value = "3.743" }
    @{ Ref = "E34"; Value = "  +2.38%  " }
    @{ Ref = "E35"; Value = "  -4.01%  " }
    @{ Ref = "D36"; Value = "5.261" }
    @{ Ref = "E36"; Value = "  +0.69%  " }
    @{ Ref = "D37"; Value = "0.06332" }
    @{ Ref = "E37"; Value = "  -4.41%  " }
    @{ Ref = "D38"; Value = "0.2143" }
    @{ Ref = "E38"; Value = "  -3.35%  " }
    @{ Ref = "D39"; Value = "1.172" }
    @{ Ref = "E39"; Value = "  -5.68%  " }
    @{ Ref = "E40"; Value = "  -5.82%  " }
    @{ Ref = "D41"; Value = "0.6305" }
    @{ Ref = "E41"; Value = "  -4.19%  " }
    @{ Ref = "D42"; Value = "11.23" }
    @{ Ref = "E42"; Value = "  -3.61%  " }
    @{ Ref = "D43"; Value = "1.200" }
    @{ Ref = "E43"; Value = "  -4.18%  " }
    @{ Ref = "D44"; Value = "1.000" }
    @{ Ref = "E44"; Value = "  +0.08%  " }
    @{ Ref = "D45"; Value = "0.5904" }
    @{ Ref = "E45"; Value = "  -4.51%  " }
    @{ Ref = "D46"; Value = "12.88" }
    @{ Ref = "E46"; Value = "  -3.03%  " }
    @{ Ref = "D47"; Value = "3.635" }
    @{ Ref = "E47"; Value = "  -3.36%  " }
    @{ Ref = "D48"; Value = "1.996" }
    @{ Ref = "E48"; Value = "  -3.66%  " }
    @{ Ref = "D49"; Value = "122.61" }
    @{ Ref = "E49"; Value = "  -2.66%  " }
    @{ Ref = "D50"; Value = "1.202" }
    @{ Ref = "E50"; Value = "  -3.58%  " }
    @{ Ref = "D51"; Value = "1.130" }
    @{ Ref = "E51"; Value = "  -2.92%  " }
)

foreach ($item in $cells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}
